$wb = $excel.ActiveWorkbook

# Column F ("想去人数") updates, applied identically to both the "展览"
# and "全部类型" sheets.
$updates = @{
    3  = 3128
    20 = 16
    21 = 47
    24 = 183
    26 = 23
    28 = 102
    29 = 2104
    33 = 204
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
